$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N ("Late"), which pushes the
# existing Late / Outstanding(heading) / Outstanding columns one position
# to the right (N->O, O->P, P->Q) and extends the used range to column Q.
$ws.Columns.Item(14).Insert()

# The newly inserted column should render at the same width as its
# neighbouring "In Advance" column (M).
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Make the "Repayment schedule" sheet the active tab/sheet and update its
# selected cell.
$ws.Activate()
$ws.Range("S7").Select()
